$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.425.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3728"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.98%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.09"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.52%  "

$ws.Range("E9").Value = "  -1.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07577"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.994"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.952"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.570.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("E17").Value = "  -1.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06748"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.285"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.27%  "

$ws.Range("E22").Value = "  -3.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.429.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.342"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.668"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.006"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.71%  "

$ws.Range("E30").Value = "  -0.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.748.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.049"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.160"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.24%  "

$ws.Range("E34").Value = "  -1.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.817"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08393"
$ws.Range("D36").Style = "Normal"

$ws.Range("E37").Value = "  +3.91%  "

$ws.Range("E38").Value = "  -2.67%  "

$ws.Range("E39").Value = "  -1.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06517"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.463"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6219"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.15%  "

$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.812"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5799"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.072"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.218"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.85%  "

$ws.Range("E51").Value = "  -0.15%  "
